$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ENW007 (row 8): append the new Jira id OPQA-1793 to the existing list ---
$ws.Range("B8").Value = "OPQA-3290||OPQA-3297||OPQA-3299||OPQA-3301||OPQA-3312||OPQA-3313||OPQA-3317||OPQA-3318||OPQA-3319||OPQA-1793"
$ws.Rows.Item(8).RowHeight = 75

# --- New test case ENW036 for OPQA-3295, appended as row 40 ---
# Copy formatting (borders/fill/wrap) down from the row above first,
# then overwrite with the new row's actual values (A, then C, then B,
# then D - matching the order the strings were originally authored in).
$ws.Range("A39:E39").Copy($ws.Range("A40:E40"))
$ws.Range("A40").Value = "ENW036"
$ws.Range("C40").Value = "Verify that field ""Group Authors"" in Neon should be displayed as ""Author"" label in Endnote after exporting the article from Neon to ENW."
$ws.Range("B40").Value = "OPQA-3295"
$ws.Range("D40").Value = "Y"

# --- Scroll / selection state saved with the workbook ---
$excel.ActiveWindow.ScrollRow = 36
$ws.Range("B45").Select()
